$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New date column F gets a slightly wider default column width (matches the
# auto-width Excel applies once a date format lands in a previously-empty
# column) - set it up front so it doesn't disturb the per-cell styles below.
$ws.Columns.Item(6).ColumnWidth = 11.42578125

# --- Header row (row 1): new "Offer" (bool) and "OfferEnd" (date) columns ---

# E1 "Offer" header - reuse the same header look as the other headers (copy D1's
# formatting, which carries the "Gut" style + text format) then overwrite the text.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Offer"

# F1 "OfferEnd" header - same header look, but with a date number format.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "OfferEnd"
$ws.Range("F1").NumberFormat = "mm-dd-yy"

# --- Data rows: Offer (bool) values ---
$ws.Range("E2").Value = $false
$ws.Range("E3").Value = $true
$ws.Range("E5").Value = $false

# --- Data rows: OfferEnd (date) value ---
# Write the raw date serial (2015-12-31) directly, then apply the short-date
# number format - assigning a real DateTime object here makes the engine
# stamp an implicit "m/d/yyyy" custom format first, leaving an orphaned
# numFmt entry behind.
$ws.Range("F3").Value = 42369
$ws.Range("F3").NumberFormat = "mm-dd-yy"

# Match the author's final selection/cursor position.
[void]$ws.Range("E5").Select()
